$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 40: update existing Start/End time entries
$ws.Range("B40").Value = 0.75
$ws.Range("C40").Value = 0.8125

# Row 43: fill in the previously-empty log entry
$ws.Range("A43").Value = 43069
$ws.Range("B43").Value = 0.70833333333333337
$ws.Range("C43").Value = 0.76388888888888884

# Update the active selection to C44
$ws.Range("C44").Select()
